# The experiment's image list was trimmed: several stimulus rows were
# removed (row 9 "images/westminster_abbey_a.jpg" and the trailing rows
# 13-16 "images/wrigley_field_b.jpg", "images/beetle.png",
# "images/john_stamos.jpg", "images/bird_parrot.png"), and the active
# selection moved to J14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the single stray row (old row 9: images/westminster_abbey_a.jpg)
$ws.Rows("9").Delete()

# Remove the trailing block of rows (old rows 13-16, now shifted to 12-15
# after the deletion above: images/wrigley_field_b.jpg, images/beetle.png,
# images/john_stamos.jpg, images/bird_parrot.png)
$ws.Rows("12:15").Delete()

# Update the selected/active cell to match the saved view state
$null = $ws.Range("J14").Select()
